# bugfixes to log_reg in R
# Updated OR / OR lCI / OR uCI values (columns B:D, rows 2-16) on the
# "comb" sheet after re-running the logistic regression in R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comb")

$ws.Range("B2").Value = 1.1309312901102899
$ws.Range("C2").Value = 1.0784326442360199
$ws.Range("D2").Value = 1.18607793415198

$ws.Range("B3").Value = 1.06971455992771
$ws.Range("C3").Value = 1.0058477673878301
$ws.Range("D3").Value = 1.1379590258515

$ws.Range("B4").Value = 0.89193697823272999
$ws.Range("C4").Value = 0.70935976294855496
$ws.Range("D4").Value = 1.1142677914804699

$ws.Range("B5").Value = 1.59311166392922
$ws.Range("C5").Value = 0.65366884681395299
$ws.Range("D5").Value = 3.61111377299695

$ws.Range("B6").Value = 1.2683009701883401
$ws.Range("C6").Value = 1.17149885738621
$ws.Range("D6").Value = 1.37305232294473

$ws.Range("B7").Value = 0.66265520206337003
$ws.Range("C7").Value = 0.61595434030856799
$ws.Range("D7").Value = 0.71325985992333396

$ws.Range("B8").Value = 0.61775214777807597
$ws.Range("C8").Value = 0.54252851658903301
$ws.Range("D8").Value = 0.70487035467919001

$ws.Range("B9").Value = 0.98900302050492594
$ws.Range("C9").Value = 0.83191932389569201
$ws.Range("D9").Value = 1.17855685649418

$ws.Range("B10").Value = 0.98128666703356404
$ws.Range("C10").Value = 0.63917838094839496
$ws.Range("D10").Value = 1.5046406396865599

$ws.Range("B11").Value = 0.60727253559376304
$ws.Range("C11").Value = 0.54724986293529398
$ws.Range("D11").Value = 0.67449379385817099

$ws.Range("B12").Value = 1.2111381788053199
$ws.Range("C12").Value = 1.1541079907199501
$ws.Range("D12").Value = 1.2710907051537099

$ws.Range("B13").Value = 1.19097545971218
$ws.Range("C13").Value = 1.1191303678253299
$ws.Range("D13").Value = 1.2678394873944101

$ws.Range("B14").Value = 1.1381097525930399
$ws.Range("C14").Value = 0.874616621820826
$ws.Range("D14").Value = 1.46920329728302

$ws.Range("B15").Value = 0.63663881685109602
$ws.Range("C15").Value = 0.0303144684844177
$ws.Range("D15").Value = 4.6221487017009402

$ws.Range("B16").Value = 1.24999772139442
$ws.Range("C16").Value = 1.1540615675438199
$ws.Range("D16").Value = 1.35374559920814

# Match the saved view state: "comb" becomes the active sheet/tab, with
# cell F27 selected (previously "proc_comb" was the active tab).
$ws.Activate()
$ws.Range("F27").Select()
